# Update cryptos list - Thu Sep 14 13:50:01 UTC 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so values such as
# "1.00", "0.410" or "213.35" are not coerced into numbers/dates, which
# would drop trailing zeros or otherwise reformat the text.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.683.14"
$ws.Range("E2").Value = "  +1.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.628.94"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.35"
$ws.Range("E5").Value = "  +0.55%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.17%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.488"
$ws.Range("E7").Value = "  +0.76%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.86%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.95%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.02"
$ws.Range("E10").Value = "  +4.89%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0832"
$ws.Range("E11").Value = "  +2.32%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.855.90"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.610.07"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.49%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +2.17%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.655.86"
$ws.Range("E16").Value = "  +1.67%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.05"
$ws.Range("E17").Value = "  +2.84%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  +0.47%  "

# Row 19 - was BitcoinCash, now Dai
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  -0.16%  "

# Row 20 - was Dai, now BitcoinCash
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "208.35"
$ws.Range("E20").Value = "  +2.89%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "4.32"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.43"
$ws.Range("E22").Value = "  +1.83%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +1.77%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.45%  "

# Row 25 - Monero
$ws.Range("D25").Value = "145.58"
$ws.Range("E25").Value = "  +1.02%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.22%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.83%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "15.40"
$ws.Range("E28").Value = "  +1.26%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "6.65"
$ws.Range("E29").Value = "  +1.69%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +7.01%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.61%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +1.79%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +0.92%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +1.92%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.42"
$ws.Range("E35").Value = "  -0.31%  "

# Row 36 - Maker
$ws.Range("D36").Value = "1.164.93"
$ws.Range("E36").Value = "  +0.54%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -0.43%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "0.810"
$ws.Range("E38").Value = "  +2.21%  "

# Row 39 - PaxDollar
$ws.Range("E39").Value = "  -0.17%  "

# Row 40 - MXToken
$ws.Range("E40").Value = "  -0.14%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +0.70%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +3.56%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "0.789"
$ws.Range("E43").Value = "  +1.14%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "1.765.98"
$ws.Range("E44").Value = "  +1.66%  "

# Row 45 - Quant
$ws.Range("D45").Value = "92.69"
$ws.Range("E45").Value = "  +0.88%  "

# Row 46 - RenderToken
$ws.Range("E46").Value = "  +2.59%  "

# Row 47 - Aave
$ws.Range("D47").Value = "54.67"
$ws.Range("E47").Value = "  +1.16%  "

# Row 48 - was Cronos, now BabyDogeCoin
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0101"
$ws.Range("E48").Value = "  +5.47%  "

# Row 49 - was Mantle, now Cronos
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  +1.01%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "7.57"
$ws.Range("E50").Value = "  +5.37%  "

# Row 51 - was USDD, now Mantle
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.410"
$ws.Range("E51").Value = "  +0.76%  "
